# Auto-generated edit script: refreshes the crypto price/volume table
# (commit: "Updated cryptos list on Wed Oct 25 23:35:30 UTC 2023 with
# GitHub Actions"). Updates Price (D) and Volume(1h) (E) columns for the
# whole list, and additionally fixes rows 15/16 where Polygon and
# WrappedBTC (name, link, price, volume) had been swapped.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '34.514.10'
$ws.Range("E2").Value = '  +1.63%  '
$ws.Range("D3").Value = '1.787.69'
$ws.Range("E3").Value = '  -0.01%  '
$ws.Range("E4").Value = '  -0.15%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '222.40'
$ws.Range("E5").Value = '  -1.66%  '
$ws.Range("E6").Value = '  -1.08%  '
$ws.Range("E7").Value = '  -0.19%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '32.32'
$ws.Range("E8").Value = '  +6.84%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.280'
$ws.Range("E9").Value = '  +0.26%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0685'
$ws.Range("E10").Value = '  +2.75%  '
$ws.Range("E11").Value = '  +1.13%  '
$ws.Range("D12").Value = '2.044.15'
$ws.Range("E12").Value = '  -0.03%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '10.99'
$ws.Range("E13").Value = '  +5.71%  '
$ws.Range("D14").Value = '1.789.49'
$ws.Range("E14").Value = '  +0.11%  '
$ws.Range("B15").Value = 'WrappedBTC'
$ws.Range("C15").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D15").Value = '34.508.00'
$ws.Range("E15").Value = '  +1.51%  '
$ws.Range("B16").Value = 'Polygon'
$ws.Range("C16").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.630'
$ws.Range("E16").Value = '  +0.80%  '
$ws.Range("E17").Value = '  +2.00%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '68.66'
$ws.Range("E18").Value = '  -0.64%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '253.90'
$ws.Range("E19").Value = '  +0.65%  '
$ws.Range("D20").Value = '0.0₃0780'
$ws.Range("E20").Value = '  +5.35%  '
$ws.Range("E21").Value = '  -0.03%  '
$ws.Range("E22").Value = '  +1.46%  '
$ws.Range("E23").Value = '  -1.35%  '
$ws.Range("E24").Value = '  +0.09%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '160.39'
$ws.Range("E25").Value = '  +1.30%  '
$ws.Range("E26").Value = '  -0.77%  '
$ws.Range("E27").Value = '  +1.22%  '
$ws.Range("E28").Value = '  -0.66%  '
$ws.Range("E29").Value = '  -0.21%  '
$ws.Range("E30").Value = '  +0.16%  '
$ws.Range("E31").Value = '  -2.05%  '
$ws.Range("E32").Value = '  -0.54%  '
$ws.Range("E33").Value = '  -0.88%  '
$ws.Range("E34").Value = '  +0.73%  '
$ws.Range("D35").Value = '1.434.89'
$ws.Range("E35").Value = '  -4.58%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.639'
$ws.Range("E36").Value = '  +0.80%  '
$ws.Range("E37").Value = '  -1.10%  '
$ws.Range("E38").Value = '  +2.30%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '85.07'
$ws.Range("E39").Value = '  +1.93%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.79'
$ws.Range("E40").Value = '  +2.89%  '
$ws.Range("E41").Value = '  -0.50%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.917'
$ws.Range("E42").Value = '  +1.74%  '
$ws.Range("E43").Value = '  +2.25%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '5.98'
$ws.Range("E44").Value = '  +4.30%  '
$ws.Range("E45").Value = '  -1.30%  '
$ws.Range("E46").Value = '  -5.49%  '
$ws.Range("D47").Value = '1.943.33'
$ws.Range("E47").Value = '  +0.08%  '
$ws.Range("E48").Value = '  +2.10%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '103.93'
$ws.Range("E49").Value = '  +6.07%  '
$ws.Range("E50").Value = '  -0.19%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '49.91'
$ws.Range("E51").Value = '  -2.85%  '
